$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting all existing data rows (23-139) down by one
# (new row N = old row N-1 for N = 24..140), and populate the newly inserted row 23
# with the new record added in this revision.
$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 45243
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 300000001
$ws.Cells.Item(23, 7).Value = "Rabanito"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 35
$ws.Cells.Item(23, 11).Value = 9000
$ws.Cells.Item(23, 12).Value = 9000
$ws.Cells.Item(23, 13).Value = 9000
$ws.Cells.Item(23, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(23, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(23, 16).Value = 750
$ws.Cells.Item(23, 17).Value = 12
$ws.Cells.Item(23, 18).Value = "Hortaliza"
